$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# zh-cn sheet: row2 gets the new handback run's timestamps,
# row3 gets the timestamps that previously belonged to row2.
$wsZhCn.Range("E2").Value = "2016-03-21 04:41:03"
$wsZhCn.Range("H2").Value = "2016-03-21 04:41:58"
$wsZhCn.Range("E3").Value = "2016-03-21 04:38:37"
$wsZhCn.Range("H3").Value = "2016-03-21 04:39:41"

# de-de sheet: same pattern.
$wsDeDe.Range("E2").Value = "2016-03-21 04:41:12"
$wsDeDe.Range("H2").Value = "2016-03-21 04:42:14"
$wsDeDe.Range("E3").Value = "2016-03-21 04:38:50"
$wsDeDe.Range("H3").Value = "2016-03-21 04:39:56"
